$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("format1")

# Core data edit: C6 previously held a blank/space string; now holds a numeric value.
$ws.Range("C6").Value = 3000.1

# Reflect the resulting active cell/selection (the cell that was edited).
$ws.Range("C7").Select()

$wb.Save()
